$d = $word.ActiveDocument

# 1) Update the title text.
$d.Content.Find.Execute(
    "數學 - 應用題 - 典型應用題 - 燕尾定理 - 一外一內比2", $true, $false, $false, $false, $false,
    $true, 1, $false, "數學 - 應用題 - 典型應用題 - 燕尾定理 - 話語霸權", 2)

# 2) Replace the text of question (1) and drop its trailing line break / image.
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = "(1) 21世紀資本論"

# 3) Remove everything after question (1): the image, and questions (2)-(4)
#    together with their images, leaving only the title and question (1).
if ($d.Paragraphs.Count -gt 2) {
    $start = $d.Paragraphs.Item(3).Range.Start
    $end = $d.Paragraphs.Item($d.Paragraphs.Count).Range.End
    $r = $d.Range($start, $end)
    $r.Delete()
}
